# link drive alla milestone
#
# 1) Remove the "justify" paragraph alignment from the three
#    "Per comandare wiCanon..." intro paragraphs (keeps the title
#    "Descrizione" paragraph's justification untouched).
# 2) Delete the trailing empty paragraph and the "Ore di lavoro: ..."
#    paragraph at the end of the document.
#
# (Run merges visible in the diff collapse identically-formatted runs
#  into one run with the same concatenated text; they produce no
#  observable change to the Word object model / rendered document, so
#  no action is required for those here.)

$d = $word.ActiveDocument

# --- 1. Clear "Justify" alignment on the three affected paragraphs ---
# "wiCanon è un programma interattivo..."
$d.Paragraphs.Item(4).Alignment = 0
# "Avviare wiCanon PC aprendo wiCanon.exe"
$d.Paragraphs.Item(5).Alignment = 0
# "Per comandare wiCanon basterà installare..."
$d.Paragraphs.Item(6).Alignment = 0

# --- 2. Remove the last two paragraphs (empty paragraph + "Ore di lavoro: ...") ---
# Deleting paragraph 19's range (the empty paragraph) merges its mark away,
# pulling the "Ore di lavoro" paragraph's text up into position 19.
$d.Paragraphs.Item(19).Range.Delete()
# Deleting the (now) paragraph 19's range removes the "Ore di lavoro" text
# and paragraph mark entirely, leaving "Nuove funzionalità del telecomando ?"
# as the final paragraph before the section break.
$d.Paragraphs.Item(19).Range.Delete()
